$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) INCM_RVN_ENUM: insert a "Reset" row right below the header (new row 2),
#    pushing R1..R5 down to rows 3..7.
# ---------------------------------------------------------------------------
$rvn = $wb.Worksheets.Item("INCM_RVN_ENUM")
$rvn.Rows("2:2").Insert()
$rvn.Range("A3:C3").Copy()
$rvn.Range("A2:C2").PasteSpecial(-4122)   # xlPasteFormats
$rvn.Range("A2").Value = "-"
$rvn.Range("B2").Value = "Reset"
$rvn.Range("C2").Value = "Reset value"
$rvn.Rows("2:2").RowHeight = 36
$rvn.AutoFilterMode = $false
$rvn.Range("A1:C7").AutoFilter()
$wb.Names.Item("INCM_RVN_ENUM!_FilterDatabase").RefersTo = "='INCM_RVN_ENUM'!`$A`$1:`$C`$7"

# ---------------------------------------------------------------------------
# 2) GEN_FRQNCY_ENUM: insert a "Reset" row right below the header (new row 2),
#    pushing the 34 existing rows down by one (A..-_Z become rows 3..36).
# ---------------------------------------------------------------------------
$frq = $wb.Worksheets.Item("GEN_FRQNCY_ENUM")
$frq.Rows("2:2").Insert()
$frq.Range("A3:C3").Copy()
$frq.Range("A2:C2").PasteSpecial(-4122)   # xlPasteFormats
$frq.Range("A2").Value = "-"
$frq.Range("B2").Value = "Reset"
$frq.Range("C2").Value = "Reset value"
$frq.Rows("2:2").RowHeight = 36
$frq.AutoFilterMode = $false
$frq.Range("A1:C36").AutoFilter()
$wb.Names.Item("GEN_FRQNCY_ENUM!_FilterDatabase").RefersTo = "='GEN_FRQNCY_ENUM'!`$A`$1:`$C`$36"

# ---------------------------------------------------------------------------
# 3) DIMS: STR32_ID -> ID (IID type column), and the EXPENSE/TP row's enum
#    reference INCM_RVN_ENUM -> INCM_EXPNS_ENUM (the REVENUE/TP row keeps
#    pointing at INCM_RVN_ENUM).
# ---------------------------------------------------------------------------
$dims = $wb.Worksheets.Item("DIMS")
$dims.Range("C2").Value = "ID"
$dims.Range("C6").Value = "ID"
$dims.Range("C7").Value = "INCM_EXPNS_ENUM"

# ---------------------------------------------------------------------------
# 4) New sheet INCM_EXPNS_ENUM (expense-type enumeration), placed after
#    ATTR_CNFDNTLTY_STTS_ENUM. Built by duplicating INCM_RVN_ENUM's layout
#    (keeps header hyperlinks / style banding / frozen header pane) and then
#    overwriting the data rows with the expense codes E1..E10 (plus Reset).
# ---------------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("INCM_RVN_ENUM")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $lastSheet)
$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "INCM_EXPNS_ENUM"

# Column widths differ slightly from the source sheet.
$newWs.Columns("A").ColumnWidth = 6.5
$newWs.Columns("B").ColumnWidth = 50
$newWs.Columns("C").ColumnWidth = 85

# The copied sheet has rows 2..6 with the correct banded styles already
# (row2/4/6 = "even" band, row3/5 = "odd" band). Grow it to 12 rows (1
# header + 11 data rows) by writing placeholder values into rows 7..12 and
# then pasting the banded formatting down from rows 6/5 alternately so the
# stripe pattern keeps going without minting brand-new style indices.
for ($r = 7; $r -le 12; $r++) {
    $newWs.Range("A$r").Value = "x"
}
for ($r = 7; $r -le 12; $r++) {
    if ((($r - 2) % 2) -eq 0) {
        $tmpl = "A6:C6"   # "even" band template (style used on row2/4/6)
    } else {
        $tmpl = "A5:C5"   # "odd" band template (style used on row3/5)
    }
    $newWs.Range($tmpl).Copy()
    $newWs.Range("A$r`:C$r").PasteSpecial(-4122)   # xlPasteFormats
    $newWs.Rows("$r`:$r").RowHeight = 36
}

$expenseRows = @(
    @("-",  "Reset",                                                                              "Reset value"),
    @("E1", "Interest paid",                                                                       "Interest paid"),
    @("E2", "Rents paid",                                                                           "Rents paid"),
    @("E3", "Financial fees/services paid by the fund (mostly for buying and selling financial assets)", "Financial fees/services paid by the fund (mostly for buying and selling financial assets)"),
    @("E4", "Management fees",                                                                      "Management fees"),
    @("E5", "Depository fees",                                                                      "Depository fees"),
    @("E6", "Taxes paid (but not those of behalf of shareholders)",                                 "Taxes paid (but not those of behalf of shareholders)"),
    @("E7", "wages",                                                                                "wages"),
    @("E8", "bad debt expense",                                                                     "bad debt expense"),
    @("E9", "bad loan expense",                                                                     "bad loan expense"),
    @("E10","Other fees and charges paid",                                                          "Other fees and charges paid")
)

for ($i = 0; $i -lt $expenseRows.Length; $i++) {
    $r = 2 + $i
    $newWs.Range("A$r").Value = $expenseRows[$i][0]
    $newWs.Range("B$r").Value = $expenseRows[$i][1]
    $newWs.Range("C$r").Value = $expenseRows[$i][2]
}

$newWs.AutoFilterMode = $false
$newWs.Range("A1:C12").AutoFilter()

# Move the selection/frozen-pane cursor back to A1 like the other sheets.
$newWs.Range("A1").Select()
